# Fan duty cycle added to gui
# Updates the "Translation" sheet of the TouchGFX texts workbook:
#   - F16 ("temperature" row): translation text changes from
#     "Temperature" to "Temp."
#   - Adds 8 new rows (20-27) for the newly added TouchGFX text ids,
#     including the new "fanDutyCycle" / "fanDutyCycleUnit" texts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Shorten the existing "temperature" translation text.
$ws.Range("F16").Value = "Temp."

# New rows appended after the existing data (row 19 was the last one).
$ws.Range("B20").Value = "SingleUseId11"
$ws.Range("C20").Value = "Default"
$ws.Range("D20").Value = "Left"
$ws.Range("E20").Value = "LTR"
$ws.Range("F20").Value = "V"

$ws.Range("B21").Value = "SingleUseId12"
$ws.Range("C21").Value = "Small"
$ws.Range("D21").Value = "Left"
$ws.Range("E21").Value = "LTR"
$ws.Range("F21").Value = "Label"

$ws.Range("B22").Value = "SingleUseId13"
$ws.Range("C22").Value = "LCD_Default"
$ws.Range("D22").Value = "Right"
$ws.Range("E22").Value = "LTR"
$ws.Range("F22").Value = "<value>"

$ws.Range("B23").Value = "fanDutyCycle"
$ws.Range("C23").Value = "Small"
$ws.Range("D23").Value = "Left"
$ws.Range("E23").Value = "LTR"
$ws.Range("F23").Value = "Fan"

$ws.Range("B24").Value = "fanDutyCycleUnit"
$ws.Range("C24").Value = "Default"
$ws.Range("D24").Value = "Left"
$ws.Range("E24").Value = "LTR"
$ws.Range("F24").Value = "%"

$ws.Range("B25").Value = "SingleUseId14"
$ws.Range("C25").Value = "Default"
$ws.Range("D25").Value = "Left"
$ws.Range("E25").Value = "LTR"
$ws.Range("F25").Value = "V"

$ws.Range("B26").Value = "SingleUseId15"
$ws.Range("C26").Value = "Small"
$ws.Range("D26").Value = "Left"
$ws.Range("E26").Value = "LTR"
$ws.Range("F26").Value = "Label"

$ws.Range("B27").Value = "SingleUseId16"
$ws.Range("C27").Value = "LCD_Default"
$ws.Range("D27").Value = "Right"
$ws.Range("E27").Value = "LTR"
$ws.Range("F27").Value = "<value>"
